$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# hole_id values for rows 2..49 (train set index -> borehole id), in row order
$holeIds = @(
    "LBU_05_08","LBU_05_22","LBU_05_07","LBU_98_2","LBU_05_12","LBU_05_13","LBU_05_29","LBU_05_30",
    "LBU_05_11","MHZ_12_03","LBU_05_14","LBU_87_6","LBU_05_24","LBU_87_4","MHZ_08_05","LBU_05_09",
    "LBU_05_06","LBU_05_15","LBU_07_01","LBU_05_28","MHZ_08_03","LBU_07_03","MHZ_08_02","LBU_98_1",
    "LBU_01_3","LBU_98_7","LBU_96_4","LBU_05_17","LBU_01_1","LBU_02_3","LBU_87_5","LBU_87_1",
    "LBU_96_1","LBU_07_02","MHZ_12_04","LBU_05_26","LBU_96_2","LBU_87_2","LBU_05_05","LBU_98_6",
    "LBU_05_04","LBU_05_19","LBU_87_3","LBU_01_2","LBU_05_01","MHZ_08_04","LBU_96_3","MHZ_12_02"
)

# New header cell A1 ("hole_id"). Column A previously had no header here, so
# borrow the bold/centered/bordered header formatting already used by B1:M1.
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Replace the old numeric 0..47 index in column A with the real hole_id text,
# keeping the existing bold/centered/bordered cell style already on A2:A49.
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
